$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "38+46=84"
$t.Rows.Item(1).Cells.Item(2).Range.Text = "56+32=88"
$t.Rows.Item(1).Cells.Item(3).Range.Text = "40-19=21"
$t.Rows.Item(1).Cells.Item(4).Range.Text = "38+8=46"
$t.Rows.Item(1).Cells.Item(5).Range.Text = "73-22=51"

$t.Rows.Item(2).Cells.Item(1).Range.Text = "23-2=21"
$t.Rows.Item(2).Cells.Item(2).Range.Text = "58-26=32"
$t.Rows.Item(2).Cells.Item(3).Range.Text = "13-6=7"
$t.Rows.Item(2).Cells.Item(4).Range.Text = "99-13=86"
$t.Rows.Item(2).Cells.Item(5).Range.Text = "32-6=26"

$t.Rows.Item(3).Cells.Item(1).Range.Text = "80-54=26"
$t.Rows.Item(3).Cells.Item(2).Range.Text = "37-30=7"
$t.Rows.Item(3).Cells.Item(3).Range.Text = "84+9=93"
$t.Rows.Item(3).Cells.Item(4).Range.Text = "32+42=74"
$t.Rows.Item(3).Cells.Item(5).Range.Text = "49-31=18"

$t.Rows.Item(4).Cells.Item(1).Range.Text = "30+22=52"
$t.Rows.Item(4).Cells.Item(2).Range.Text = "8+47=55"
$t.Rows.Item(4).Cells.Item(3).Range.Text = "87-47=40"
$t.Rows.Item(4).Cells.Item(4).Range.Text = "26+15=41"
$t.Rows.Item(4).Cells.Item(5).Range.Text = "94-13=81"

$t.Rows.Item(5).Cells.Item(1).Range.Text = "82-40=42"
$t.Rows.Item(5).Cells.Item(2).Range.Text = "16+71=87"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "82-73=9"
$t.Rows.Item(5).Cells.Item(4).Range.Text = "2+24=26"
$t.Rows.Item(5).Cells.Item(5).Range.Text = "9+82=91"

$t.Rows.Item(6).Cells.Item(1).Range.Text = "52+11=63"
$t.Rows.Item(6).Cells.Item(2).Range.Text = "82-27=55"
$t.Rows.Item(6).Cells.Item(3).Range.Text = "2+0=2"
$t.Rows.Item(6).Cells.Item(4).Range.Text = "63-29=34"
$t.Rows.Item(6).Cells.Item(5).Range.Text = "91-58=33"

$t.Rows.Item(7).Cells.Item(1).Range.Text = "78-66=12"
$t.Rows.Item(7).Cells.Item(2).Range.Text = "57-23=34"
$t.Rows.Item(7).Cells.Item(3).Range.Text = "29-23=6"
$t.Rows.Item(7).Cells.Item(4).Range.Text = "38+31=69"
$t.Rows.Item(7).Cells.Item(5).Range.Text = "25+18=43"

$t.Rows.Item(8).Cells.Item(1).Range.Text = "2+56=58"
$t.Rows.Item(8).Cells.Item(2).Range.Text = "5-4=1"
$t.Rows.Item(8).Cells.Item(3).Range.Text = "9+34=43"
$t.Rows.Item(8).Cells.Item(4).Range.Text = "96-32=64"
$t.Rows.Item(8).Cells.Item(5).Range.Text = "50-7=43"

$t.Rows.Item(9).Cells.Item(1).Range.Text = "48+39=87"
$t.Rows.Item(9).Cells.Item(2).Range.Text = "14+0=14"
$t.Rows.Item(9).Cells.Item(3).Range.Text = "47-8=39"
$t.Rows.Item(9).Cells.Item(4).Range.Text = "13-9=4"
$t.Rows.Item(9).Cells.Item(5).Range.Text = "4-2=2"

$t.Rows.Item(10).Cells.Item(1).Range.Text = "1+53=54"
$t.Rows.Item(10).Cells.Item(2).Range.Text = "7+31=38"
$t.Rows.Item(10).Cells.Item(3).Range.Text = "13+35=48"
$t.Rows.Item(10).Cells.Item(4).Range.Text = "48-31=17"
$t.Rows.Item(10).Cells.Item(5).Range.Text = "38+4=42"

$t.Rows.Item(11).Cells.Item(1).Range.Text = "18+15=33"
$t.Rows.Item(11).Cells.Item(2).Range.Text = "69+17=86"
$t.Rows.Item(11).Cells.Item(3).Range.Text = "15+72=87"
$t.Rows.Item(11).Cells.Item(4).Range.Text = "48+51=99"
$t.Rows.Item(11).Cells.Item(5).Range.Text = "94-33=61"

$t.Rows.Item(12).Cells.Item(1).Range.Text = "45-31=14"
$t.Rows.Item(12).Cells.Item(2).Range.Text = "54-13=41"
$t.Rows.Item(12).Cells.Item(3).Range.Text = "59+35=94"
$t.Rows.Item(12).Cells.Item(4).Range.Text = "70-47=23"
$t.Rows.Item(12).Cells.Item(5).Range.Text = "57+9=66"

$t.Rows.Item(13).Cells.Item(1).Range.Text = "73+0=73"
$t.Rows.Item(13).Cells.Item(2).Range.Text = "70-31=39"
$t.Rows.Item(13).Cells.Item(3).Range.Text = "5+69=74"
$t.Rows.Item(13).Cells.Item(4).Range.Text = "3+91=94"
$t.Rows.Item(13).Cells.Item(5).Range.Text = "31-13=18"

$t.Rows.Item(14).Cells.Item(1).Range.Text = "55-52=3"
$t.Rows.Item(14).Cells.Item(2).Range.Text = "96-44=52"
$t.Rows.Item(14).Cells.Item(3).Range.Text = "20+31=51"
$t.Rows.Item(14).Cells.Item(4).Range.Text = "89-46=43"
$t.Rows.Item(14).Cells.Item(5).Range.Text = "36+2=38"

$t.Rows.Item(15).Cells.Item(1).Range.Text = "77-12=65"
$t.Rows.Item(15).Cells.Item(2).Range.Text = "79-74=5"
$t.Rows.Item(15).Cells.Item(3).Range.Text = "16+4=20"
$t.Rows.Item(15).Cells.Item(4).Range.Text = "22-5=17"
$t.Rows.Item(15).Cells.Item(5).Range.Text = "14+7=21"

$t.Rows.Item(16).Cells.Item(1).Range.Text = "15+83=98"
$t.Rows.Item(16).Cells.Item(2).Range.Text = "93-17=76"
$t.Rows.Item(16).Cells.Item(3).Range.Text = "65+30=95"
$t.Rows.Item(16).Cells.Item(4).Range.Text = "32+43=75"
$t.Rows.Item(16).Cells.Item(5).Range.Text = "71+16=87"

$t.Rows.Item(17).Cells.Item(1).Range.Text = "93-41=52"
$t.Rows.Item(17).Cells.Item(2).Range.Text = "6+54=60"
$t.Rows.Item(17).Cells.Item(3).Range.Text = "96-21=75"
$t.Rows.Item(17).Cells.Item(4).Range.Text = "59-5=54"
$t.Rows.Item(17).Cells.Item(5).Range.Text = "17-14=3"

$t.Rows.Item(18).Cells.Item(1).Range.Text = "96-85=11"
$t.Rows.Item(18).Cells.Item(2).Range.Text = "59-37=22"
$t.Rows.Item(18).Cells.Item(3).Range.Text = "88-12=76"
$t.Rows.Item(18).Cells.Item(4).Range.Text = "71+20=91"
$t.Rows.Item(18).Cells.Item(5).Range.Text = "69-41=28"

$t.Rows.Item(19).Cells.Item(1).Range.Text = "25-11=14"
$t.Rows.Item(19).Cells.Item(2).Range.Text = "42+24=66"
$t.Rows.Item(19).Cells.Item(3).Range.Text = "83-6=77"
$t.Rows.Item(19).Cells.Item(4).Range.Text = "37+56=93"
$t.Rows.Item(19).Cells.Item(5).Range.Text = "42-21=21"

$t.Rows.Item(20).Cells.Item(1).Range.Text = "92-82=10"
$t.Rows.Item(20).Cells.Item(2).Range.Text = "36+13=49"
$t.Rows.Item(20).Cells.Item(3).Range.Text = "24-0=24"
$t.Rows.Item(20).Cells.Item(4).Range.Text = "63-49=14"
$t.Rows.Item(20).Cells.Item(5).Range.Text = "24+52=76"
